$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("FEINmismatch")
$sheet1Dates = @(
    "Fri Sep 29 11:16:53 EDT 2023",
    "Fri Sep 29 11:17:08 EDT 2023",
    "Fri Sep 29 11:17:20 EDT 2023",
    "Fri Sep 29 11:17:32 EDT 2023",
    "Fri Sep 29 11:17:45 EDT 2023",
    "Fri Sep 29 11:17:57 EDT 2023",
    "Fri Sep 29 11:18:09 EDT 2023",
    "Fri Sep 29 11:18:20 EDT 2023",
    "Fri Sep 29 11:18:32 EDT 2023",
    "Fri Sep 29 11:18:44 EDT 2023",
    "Fri Sep 29 11:18:56 EDT 2023",
    "Fri Sep 29 11:19:08 EDT 2023"
)

for ($i = 0; $i -lt $sheet1Dates.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $sheet1Dates[$i]
}

$ws2 = $wb.Worksheets.Item("FEINSSNmismatch")
$sheet2Dates = @(
    "Fri Sep 29 11:19:21 EDT 2023",
    "Fri Sep 29 11:19:32 EDT 2023",
    "Fri Sep 29 11:19:43 EDT 2023",
    "Fri Sep 29 11:19:55 EDT 2023",
    "Fri Sep 29 11:20:06 EDT 2023",
    "Fri Sep 29 11:20:18 EDT 2023",
    "Fri Sep 29 11:20:29 EDT 2023",
    "Fri Sep 29 11:20:40 EDT 2023",
    "Fri Sep 29 11:20:52 EDT 2023",
    "Fri Sep 29 11:21:03 EDT 2023",
    "Fri Sep 29 11:21:14 EDT 2023",
    "Fri Sep 29 11:21:26 EDT 2023",
    "Fri Sep 29 11:21:37 EDT 2023",
    "Fri Sep 29 11:21:49 EDT 2023",
    "Fri Sep 29 11:22:00 EDT 2023",
    "Fri Sep 29 11:22:11 EDT 2023"
)

for ($i = 0; $i -lt $sheet2Dates.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 2).Value = $sheet2Dates[$i]
}
